# Fazendo alterações na Apresentação e no Backlog
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix capitalization of the "Response web" section header -> "Response Web"
$ws.Range("A11").Value = "Response Web"

# Match the font styling used by the other section headers (A7, A14, A16): red, size 16
$ws.Range("A11").Font.Color = 255
$ws.Range("A11").Font.Size = 16

# Reflect the resulting selection/view state (user clicked into A11 after the edit)
$ws.Range("A11").Select()
